$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column J (the previous last data column) into column K
# so the new 2023 column visually matches the rest of the table.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)

# Fill in the new 2023 data.
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 3469
$ws.Range("K5").Value = 1612
$ws.Range("K6").Value = 1857

# Give column K (now the right-most column of the table) a thin right border
# to close off the table, matching the existing top/bottom rule lines.
$ws.Range("K3:K6").Borders.Item(10).LineStyle = 1
$ws.Range("K3:K6").Borders.Item(10).Weight = 2

# Match column widths used by the rest of the data columns.
$ws.Range("K1:K6").ColumnWidth = 8.7109375
